# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" (holdings detail, same layout as the other
#    quarterly sheets) positioned right before the "总计" (totals) sheet.
# 2. Prepend a new "2022-Q1" summary row to the "总计" sheet, pushing the
#    existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "2022-Q1" worksheet
# ---------------------------------------------------------------------------

$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")

$ws = $wb.Worksheets.Add($totalSheet)
$ws.Name = "2022-Q1"

# Headers (row 1, columns B..H) - copy style from an existing quarter sheet
# so formatting (bold, centered, bordered) matches exactly.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$templateSheet.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

# Holdings data, in the exact order/values from the source report.
$rows = @(
    @("166005", "中欧价值发现混合 -A", "43.52", "93.97", "3.25", "1.4144", 8),
    @("001882", "中欧价值发现混合 -E", "43.52", "93.97", "3.25", "1.4144", 8),
    @("001810", "中欧潜力价值灵活配置混合A", "28.67", "94.05", "2.63", "0.7540", 9),
    @("010744", "工银瑞信灵动价值混合A", "15.18", "76.76", "4.63", "0.7028", 2),
    @("160916", "大成优选混合(LOF)", "16.14", "89.35", "3.55", "0.5730", 10),
    @("008269", "大成睿享混合A", "17.69", "65.25", "3.05", "0.5395", 6),
    @("004232", "中欧价值发现混合 -C", "10.98", "93.97", "3.25", "0.3568", 8),
    @("166024", "中欧恒利三年定期开放混合", "4.48", "98.71", "3.92", "0.1756", 6),
    @("001651", "工银瑞信新蓝筹股票A", "4.71", "80.99", "2.96", "0.1394", 6),
    @("481013", "工银消费服务混合A", "3.80", "61.85", "3.21", "0.1220", 5),
    @("005764", "中欧潜力价值灵活配置混合C", "3.43", "94.05", "2.63", "0.0902", 9),
    @("008270", "大成睿享混合C", "2.87", "65.25", "3.05", "0.0875", 6),
    @("010745", "工银瑞信灵动价值混合C", "1.01", "76.76", "4.63", "0.0468", 2),
    @("011476", "工银瑞信新蓝筹股票C", "0.42", "80.99", "2.96", "0.0124", 6),
    @("011475", "工银消费服务混合C", "0.03", "61.85", "3.21", "0.0010", 5)
)

$firstDataRow = 2
for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $rowNum = $firstDataRow + $r

    # Column A: running index (0-based), numeric, same style as the template.
    $ws.Cells.Item($rowNum, 1).Value = $r

    # Column B: fund code - keep as literal text even though it looks numeric.
    $ws.Cells.Item($rowNum, 2).Formula = '="' + $row[0] + '"'

    # Column C: fund name - plain text, never numeric-looking.
    $ws.Cells.Item($rowNum, 3).Value = $row[1]

    # Columns D-G: decimal-looking figures stored as literal text.
    $ws.Cells.Item($rowNum, 4).Formula = '="' + $row[2] + '"'
    $ws.Cells.Item($rowNum, 5).Formula = '="' + $row[3] + '"'
    $ws.Cells.Item($rowNum, 6).Formula = '="' + $row[4] + '"'
    $ws.Cells.Item($rowNum, 7).Formula = '="' + $row[5] + '"'

    # Column H: rank, a genuine number.
    $ws.Cells.Item($rowNum, 8).Value = $row[6]
}

$lastDataRow = $firstDataRow + $rows.Length - 1

# Convert the helper "=""text""" formulas into plain literal text values
# (matches the inline-string layout used throughout the workbook, with no
# extra number formatting / styling baked in).
$ws.Range("B$firstDataRow`:B$lastDataRow").Copy()
$ws.Range("B$firstDataRow`:B$lastDataRow").PasteSpecial(-4163)
$ws.Range("D$firstDataRow`:G$lastDataRow").Copy()
$ws.Range("D$firstDataRow`:G$lastDataRow").PasteSpecial(-4163)

# Column A numbers use the same centered/bold style as the other sheets.
$templateSheet.Range("A2").Copy()
$ws.Range("A$firstDataRow`:A$lastDataRow").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. "总计" sheet - prepend the 2022-Q1 summary row
# ---------------------------------------------------------------------------
# NOTE: the reference captured in $totalSheet before the Worksheets.Add call
# above now resolves to whatever sheet occupies that original tab position
# (i.e. the freshly inserted "2022-Q1" sheet), not "总计" itself - worksheet
# handles in this host are position-bound, not identity-bound. Re-resolve by
# name now that the sheet collection is stable again.
$totalSheet = $wb.Worksheets.Item("总计")

$summary = @(
    @("2022-Q1", 15, 6.43),
    @("2021-Q4", 18, 7.42),
    @("2021-Q3", 22, 5.01),
    @("2021-Q2", 23, 4.27),
    @("2021-Q1", 11, 4.41),
    @("2020-Q4", 9, 4.4)
)

for ($i = 0; $i -lt $summary.Length; $i++) {
    $entry = $summary[$i]
    $rowNum = $i + 2
    $totalSheet.Cells.Item($rowNum, 1).Value = $i
    $totalSheet.Cells.Item($rowNum, 2).Value = $entry[0]
    $totalSheet.Cells.Item($rowNum, 3).Value = $entry[1]
    $totalSheet.Cells.Item($rowNum, 4).Value = $entry[2]
}
